# Changed default delimiter character from '~' to '.' in header strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "elements.H"
$ws.Range("D1").Value = "elements.O"

# Update active selection to match the edit.
$ws.Range("M1").Select()
